# Refresh the cryptocurrency price/volume snapshot (coinranking.com scrape).
# Mirrors the automated GitHub Actions data-refresh commit: updates prices,
# 1h volume percentages, and re-ranks a handful of exchange-token rows whose
# relative ordering changed (GateToken/MXToken/FTXToken/WazirX/... block).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'246.81"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.06%"
$ws.Range("E2").Style = "Normal"

# Row 3
$ws.Range("D3").Value = "'29.61"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'8.16%"
$ws.Range("E3").Style = "Normal"

# Row 4
$ws.Range("D4").Value = "'5.192"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'2.62%"
$ws.Range("E4").Style = "Normal"

# Row 5
$ws.Range("D5").Value = "'0.05710"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'0.49%"
$ws.Range("E5").Style = "Normal"

# Row 6
$ws.Range("D6").Value = "'6.580"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.63%"
$ws.Range("E6").Style = "Normal"

# Row 7
$ws.Range("B7").Value = "GateToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D7").Value = "'3.070"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.29%"
$ws.Range("E7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8587"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'4.54%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.8788"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'4.57%"
$ws.Range("E9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1369"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'3.34%"
$ws.Range("E10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07072"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'2.29%"
$ws.Range("E11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02877"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'0.43%"
$ws.Range("E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09385"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.11%"
$ws.Range("E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001513"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.24%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04163"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'1.53%"
$ws.Range("E15").Style = "Normal"

# Row 16
$ws.Range("D16").Value = "'0.006190"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.90%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("E17").Value = "'5,107.59%"
$ws.Range("E17").Style = "Normal"

# Row 18
$ws.Range("D18").Value = "'3.479"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'-0.83%"
$ws.Range("E18").Style = "Normal"

# Row 19
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.282"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-1.27%"
$ws.Range("E19").Style = "Normal"

# Row 20
$ws.Range("B20").Value = "One"
$ws.Range("C20").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D20").Value = "'0.01033"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'1,622.37%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("E21").Value = "'0.75%"
$ws.Range("E21").Style = "Normal"

# Row 22
$ws.Range("D22").Value = "'0.03306"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'3.49%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("E23").Value = "'0.64%"
$ws.Range("E23").Style = "Normal"

# Row 24
$ws.Range("D24").Value = "'3.470"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'-2.48%"
$ws.Range("E24").Style = "Normal"

# Row 25
$ws.Range("D25").Value = "'0.1378"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.32%"
$ws.Range("E25").Style = "Normal"

# Row 26
$ws.Range("D26").Value = "'0.005051"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'28.03%"
$ws.Range("E26").Style = "Normal"

# Row 27
$ws.Range("D27").Value = "'0.001219"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.11%"
$ws.Range("E27").Style = "Normal"

# Row 28
$ws.Range("E28").Value = "'23.36%"
$ws.Range("E28").Style = "Normal"

# Row 40
$ws.Range("D40").Value = "'0.03756"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-0.40%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41").Value = "'0.005675"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-8.17%"
$ws.Range("E41").Style = "Normal"

# Row 42
$ws.Range("E42").Value = "'1.97%"
$ws.Range("E42").Style = "Normal"

# Row 43
$ws.Range("D43").Value = "'0.002098"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-8.77%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44").Value = "'0.009947"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'6.93%"
$ws.Range("E44").Style = "Normal"

# Row 45
$ws.Range("D45").Value = "'0.00005127"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-1.71%"
$ws.Range("E45").Style = "Normal"

# Row 47
$ws.Range("D47").Value = "'0.07094"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-30.10%"
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("D48").Value = "'0.002583"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'0.50%"
$ws.Range("E48").Style = "Normal"
